# Sync attendance_reports: normalize "Recorded By" (column G) ordering on the
# "Session Analysis Results" sheet so "System" is no longer always reported
# first (reorders the trailing two recorder entries for known combinations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, system, backup@backdoor.com" = "System, backup@backdoor.com, system";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$lastRow = $ws.UsedRange.Rows.Count()

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value()
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
